$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "hussein2492024053624"
$ws.Range("B6").Value = "Test@123"

$ws.Range("A7").Value = "hussein2492024053839"
$ws.Range("B7").Value = "Test@123"

$ws.Range("A8").Value = "hussein2492024054230"
$ws.Range("B8").Value = "Test@123"
